# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions crypto price/volume refresh on 2024-08-10)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, [string]$val)
    # Force plain-text storage so numeric-looking strings (e.g. "6.70",
    # "0.0552") are not silently coerced into Excel numbers, which would
    # drop significant trailing/leading zeros and change the cell type.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    # Restore the default (unstyled) look so we don't leave a stray
    # Text-formatted style behind on cells that were plain General before.
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = "60.260.36"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.587.39"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("E4").Value = "  +0.75%  "
Set-TextCell $ws "D5" "508.33"
$ws.Range("E5").Value = "  -0.22%  "
Set-TextCell $ws "D6" "153.17"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("E7").Value = "  +0.48%  "
Set-TextCell $ws "D8" "0.588"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").Value = "2.596.60"
$ws.Range("E9").Value = "  -2.15%  "
Set-TextCell $ws "D10" "6.70"
$ws.Range("E10").Value = "  +6.24%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").Value = "3.043.76"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "60.271.95"
$ws.Range("E15").Value = "  -0.40%  "
Set-TextCell $ws "D16" "21.45"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "2.589.78"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("E19").Value = "  -0.66%  "
Set-TextCell $ws "D20" "355.01"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("E21").Value = "  +0.25%  "
Set-TextCell $ws "D22" "6.12"
$ws.Range("E22").Value = "  -0.37%  "
Set-TextCell $ws "D23" "0.999"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("E30").Value = "  +0.20%  "
Set-TextCell $ws "D31" "19.35"
$ws.Range("E31").Value = "  -0.47%  "
Set-TextCell $ws "D32" "151.61"
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("E33").Value = "  -0.77%  "
Set-TextCell $ws "D34" "5.71"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  -2.58%  "
Set-TextCell $ws "D37" "0.865"
$ws.Range("E37").Value = "  +4.19%  "
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D39" "3.75"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws "D40" "36.08"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("E41").Value = "  -2.55%  "
Set-TextCell $ws "D42" "293.73"
$ws.Range("E42").Value = "  -3.66%  "
Set-TextCell $ws "D43" "0.101"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -4.21%  "
Set-TextCell $ws "D45" "0.998"
$ws.Range("E45").Value = "  +0.67%  "
Set-TextCell $ws "D46" "0.0552"
$ws.Range("E46").Value = "  -3.73%  "
Set-TextCell $ws "D47" "19.63"
$ws.Range("E47").Value = "  -1.15%  "
Set-TextCell $ws "D48" "4.79"
$ws.Range("E48").Value = "  -3.83%  "
Set-TextCell $ws "D49" "0.0233"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "1.986.77"
$ws.Range("E51").Value = "  -2.10%  "
